# Update "ResumoInscricoes" data: refresh Pagos (F) and Inscricoes homologadas (H)
# values for several rows, plus Inscritos (E) for a couple of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 16
$ws.Range("H2").Value = 18

# Row 3
$ws.Range("F3").Value = 16
$ws.Range("H3").Value = 19

# Row 4
$ws.Range("E4").Value = 25

# Row 6
$ws.Range("F6").Value = 26
$ws.Range("H6").Value = 33

# Row 8
$ws.Range("F8").Value = 16
$ws.Range("H8").Value = 22

# Row 11
$ws.Range("F11").Value = 16
$ws.Range("H11").Value = 17

# Row 15
$ws.Range("F15").Value = 51
$ws.Range("H15").Value = 62

# Row 16
$ws.Range("F16").Value = 94
$ws.Range("H16").Value = 181

# Row 17
$ws.Range("F17").Value = 14
$ws.Range("H17").Value = 15

# Row 18
$ws.Range("E18").Value = 99
$ws.Range("F18").Value = 33
$ws.Range("H18").Value = 56
